$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Rename the two "7m" resolution CSV inputs to "6m" everywhere they occur
# (Run-geometry column J and Validation-geometry column O) across all
# scenario group header rows (3, 11, 16, 21, ..., 91).
$rows = @(3,11,16,21,26,31,36,41,46,51,56,61,66,71,76,81,86,91)
foreach ($r in $rows) {
    $jCell = $ws.Cells.Item($r, 10)   # column J
    $oCell = $ws.Cells.Item($r, 15)   # column O
    if ($jCell.Value2 -eq "Inputs\Selwyn_XS3_7m.csv") {
        $jCell.Value = "Inputs\Selwyn_XS3_6m.csv"
    }
    if ($oCell.Value2 -eq "Inputs\Selwyn_XS7_7m.csv") {
        $oCell.Value = "Inputs\Selwyn_XS7_6m.csv"
    }
}

# Reduce the optimisation upper bound for the first scenario row (1a) so
# the 6m-resolution slope model converges.
$ws.Cells.Item(3, 13).Value = 0.5   # column M, row 3

# Move the active selection to L3 on the Scenarios sheet.
$ws.Range("L3").Select()
